# Powerpoint writer: consolidate text run nodes.
# Merge the leading "<word>" + " " runs into a single "<word> " run, leaving
# the following run (the number / next word) untouched as its own run.

$p = $ppt.ActivePresentation

# Slide 1, Title "Slide" + " " + "1" -> "Slide " + "1"
$s1 = $p.Slides.Item(1)
$title1 = $s1.Shapes.Item(1)
$title1.TextFrame.TextRange.Characters(1, 6).Text = "Slide "

# Slide 1, TextBox "an" + " " + "image" -> "an " + "image"
$tb1 = $s1.Shapes.Item(3)
$tb1.TextFrame.TextRange.Characters(1, 3).Text = "an "

# Slide 2, Title "Slide" + " " + "2" -> "Slide " + "2"
$s2 = $p.Slides.Item(2)
$title2 = $s2.Shapes.Item(1)
$title2.TextFrame.TextRange.Characters(1, 6).Text = "Slide "

# Slide 2, TextBox "an" + " " + "image" -> "an " + "image"
$tb2 = $s2.Shapes.Item(4)
$tb2.TextFrame.TextRange.Characters(1, 3).Text = "an "
